$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the B1:F1 row down to B11:F11 (copy & paste of the original table)
$ws.Range("B1:F1").Copy($ws.Range("B11:F11"))

# Copy the A2:F5 block down to A12:F15 (copy & paste of the original table)
$ws.Range("A2:F5").Copy($ws.Range("A12:F15"))

# Update the timestamp text stored in A10
$ws.Range("A10").Value = "2023-03-16 오전 8:27:04"
